# Semainiers de du 11mars
# Update the weekly tracker ("Semainier") statuses and actual hours for
# the second sprint's task list (rows 50-54 on Feuil1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 50 ("Administration des roles des utilisateurs") and row 51
# ("Pages administration et gestion par GUI") are now finished.
$ws.Range("H50").Value = "Terminé"
$ws.Range("H51").Value = "Terminé"

# Row 52 ("Pages gestion d'entreprises (Create,Edit)") - work started,
# 3 real hours logged, status moved from "En attente" to "En cours".
$ws.Range("C52").Value = 3
$ws.Range("H52").Value = "En cours"

# Row 53 ("Pages gestion d'evenement (Create,Edit,Delete)") - work
# started, 8 real hours logged, status moved from "En attente" to
# "En cours".
$ws.Range("C53").Value = 8
$ws.Range("H53").Value = "En cours"

# Row 54 ("Messageries entre utilisateur") - work started, 10 real
# hours logged, status moved from "En attente" to "En cours".
$ws.Range("F54").Value = 10
$ws.Range("H54").Value = "En cours"

# Leave the selection where the editing session ended.
$ws.Range("K41").Select()
